# Update the "About" sheet: insert a new row after row 11 to split the
# note into two lines, and adjust the note text.
$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$eci = $wb.Worksheets.Item("ECiCpCU")

# Insert a new row at row 12 (pushes old rows 12+ down by one).
$about.Rows.Item(12).Insert()

# Update row 11 text (drop trailing period) and fill in new row 12.
$about.Range("A11").Value = "For each component, we average the carbon intensities of the one or more most applicable industries"
$about.Range("A12").Value = "and divide by the component cost."
$about.Rows.Item(12).Style = $about.Range("A19").Style

# Fix the label on the ECiCpCU sheet.
$eci.Range("B1").Value = "Embedded tons CO2e/2012$"

$wb.Save()
